# Daily "remaining days" rollover update.
# Column D = total days, E = remaining days, F = start date (yyyymmdd).
# For each data row: decrement E by 1; if E was already 1 (i.e. the cycle
# finished), reset E back to the total (D) and bump the start date (F) to
# today (2026-02-05).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$today = 20260205
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $dVal = $dCell.Value2
    $eVal = $eCell.Value2
    $fVal = $fCell.Value2

    if ($dVal -eq $null -or $eVal -eq $null) {
        continue
    }

    # Skip rows whose start date is malformed (not a proper yyyymmdd,
    # i.e. not 8 digits) - they are left untouched by the refresh job.
    $fStr = [string]$fVal
    if ($fStr.Length -ne 8) {
        continue
    }

    if ($eVal -eq 1) {
        $eCell.Value2 = $dVal
        $fCell.Value2 = $today
    } else {
        $eCell.Value2 = $eVal - 1
    }
}
